# Apply the latest cryptos-list snapshot (price + 1h volume deltas) to Sheet1.
# D-column prices are free-text (e.g. "67.415.33", "173.90") so they must be
# written as Text, otherwise Excel auto-coerces them to numbers and mangles
# thousands separators / trailing zeros. Force Text format, write, then drop
# back to the "Normal" style so no stray number-format style sticks to the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $cell = $ws.Range($range)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "67.415.33"
$ws.Range("E2").Value = "  +0.38%  "

# Row 3
Set-TextValue "D3" "2.558.93"
$ws.Range("E3").Value = "  -2.29%  "

# Row 4
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
Set-TextValue "D5" "594.01"
$ws.Range("E5").Value = "  +0.82%  "

# Row 6
Set-TextValue "D6" "173.90"

# Row 7
$ws.Range("E7").Value = "  +0.06%  "

# Row 8
Set-TextValue "D8" "0.532"
$ws.Range("E8").Value = "  +0.12%  "

# Row 9
Set-TextValue "D9" "2.557.55"
$ws.Range("E9").Value = "  -2.30%  "

# Row 10
$ws.Range("E10").Value = "  +0.26%  "

# Row 11
$ws.Range("E11").Value = "  +1.81%  "

# Row 12
$ws.Range("B12").Value = "Cardano"
$ws.Range("C12").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
Set-TextValue "D12" "0.353"
$ws.Range("E12").Value = "  -3.82%  "

# Row 13
$ws.Range("B13").Value = "Toncoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D13" "5.18"
$ws.Range("E13").Value = "  -0.77%  "

# Row 14
Set-TextValue "D14" "27.16"
$ws.Range("E14").Value = "  -0.63%  "

# Row 15
Set-TextValue "D15" "3.014.59"
$ws.Range("E15").Value = "  -2.62%  "

# Row 16
$ws.Range("E16").Value = "  -0.53%  "

# Row 17
Set-TextValue "D17" "67.295.22"
$ws.Range("E17").Value = "  +0.50%  "

# Row 18
Set-TextValue "D18" "2.555.65"
$ws.Range("E18").Value = "  -2.20%  "

# Row 19
Set-TextValue "D19" "8.05"
$ws.Range("E19").Value = "  +3.26%  "

# Row 20
Set-TextValue "D20" "11.43"
$ws.Range("E20").Value = "  -2.74%  "

# Row 21
Set-TextValue "D21" "356.54"
$ws.Range("E21").Value = "  +0.00%  "

# Row 22
Set-TextValue "D22" "4.24"
$ws.Range("E22").Value = "  -1.04%  "

# Row 23
Set-TextValue "D23" "4.70"
$ws.Range("E23").Value = "  +1.07%  "

# Row 24
$ws.Range("E24").Value = "  +5.66%  "

# Row 25
$ws.Range("E25").Value = "  -0.03%  "

# Row 26
Set-TextValue "D26" "70.34"
$ws.Range("E26").Value = "  +1.37%  "

# Row 27
Set-TextValue "D27" "10.16"
$ws.Range("E27").Value = "  -3.99%  "

# Row 28
Set-TextValue "D28" "2.681.90"
$ws.Range("E28").Value = "  -2.37%  "

# Row 29
$ws.Range("E29").Value = "  -0.15%  "

# Row 30
Set-TextValue "D30" "0.0000100"
$ws.Range("E30").Value = "  +0.53%  "

# Row 31
Set-TextValue "D31" "538.43"
$ws.Range("E31").Value = "  -1.40%  "

# Row 32
Set-TextValue "D32" "8.26"
$ws.Range("E32").Value = "  +4.89%  "

# Row 33
Set-TextValue "D33" "1.39"
$ws.Range("E33").Value = "  +3.66%  "

# Row 34
$ws.Range("E34").Value = "  -0.17%  "

# Row 35
Set-TextValue "D35" "0.133"
$ws.Range("E35").Value = "  -0.01%  "

# Row 37
$ws.Range("E37").Value = "  +0.55%  "

# Row 38
Set-TextValue "D38" "157.40"
$ws.Range("E38").Value = "  -0.10%  "

# Row 39
Set-TextValue "D39" "18.82"
$ws.Range("E39").Value = "  -0.53%  "

# Row 40
$ws.Range("E40").Value = "  +1.27%  "

# Row 41
Set-TextValue "D41" "0.359"
$ws.Range("E41").Value = "  -1.79%  "

# Row 42
$ws.Range("E42").Value = "  +0.44%  "

# Row 43
Set-TextValue "D43" "5.21"
$ws.Range("E43").Value = "  +1.44%  "

# Row 44
$ws.Range("E44").Value = "  +5.03%  "

# Row 45
$ws.Range("E45").Value = "  -0.03%  "

# Row 46
Set-TextValue "D46" "39.71"
$ws.Range("E46").Value = "  -1.28%  "

# Row 47
Set-TextValue "D47" "151.48"
$ws.Range("E47").Value = "  +0.37%  "

# Row 48
Set-TextValue "D48" "0.567"
$ws.Range("E48").Value = "  -1.98%  "

# Row 49
Set-TextValue "D49" "0.0₆0282"
$ws.Range("E49").Value = "  -4.90%  "

# Row 50
Set-TextValue "D50" "3.73"
$ws.Range("E50").Value = "  -1.07%  "

# Row 51
$ws.Range("E51").Value = "  +1.15%  "
